# "mise en commentaires des paramètres de recherche qu'on n'utilise plus"
#
# Net effect on the "Metadata" sheet:
#   - a new "Jurisdiction" / "" row is inserted right after the "Contact" row,
#     pushing every row from "Description" down through "Count" down by one;
#   - the "Date" value is refreshed to the new publication timestamp.
# The "Concepts" sheet (sheet2) is untouched - its shared-string indices will
# simply be renumbered automatically because of the sharedStrings table growth.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- 1. Make room for the new row: shift rows 11..21 down into 12..22 ---

# Prime row 22 (brand new) with the same formatting (border/alignment/etc.)
# as the last existing data row, so the shift below doesn't leave it with
# the workbook's bare default style.
$ws.Range("A21:B21").Copy()
$ws.Range("A22:B22").PasteSpecial(-4122)  # xlPasteFormats

# Walk bottom-up so each source row is read before it gets overwritten.
# ClearContents first because PasteSpecial(xlPasteValues) of a blank source
# cell does not blank out a non-blank destination cell on its own.
for ($r = 21; $r -ge 11; $r--) {
    $ws.Range("A$($r+1):B$($r+1)").ClearContents()
    $ws.Range("A$r`:B$r").Copy()
    $ws.Range("A$($r+1):B$($r+1)").PasteSpecial(-4163)  # xlPasteValues
}
$excel.CutCopyMode = 0

# --- 2. Fill the freshly opened row 11 with the new "Jurisdiction" entry ---
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# --- 3. Refresh the publication Date value ---
$ws.Cells.Item(8, 2).Value = "2024-07-01T07:50:29+00:00"
